$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each target cell currently stores a text value that looks like a number
# (e.g. "308.86" or "-0.64%"). Force Text format so COM keeps it as a string
# instead of auto-converting it to a numeric/percentage value.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "308.86"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.64%"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.11%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.065"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-1.98%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07902"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.03%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.059"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "3.67%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.434"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.66%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.275"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "1.64%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-7.75%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9332"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.68%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1282"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-7.76%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1884"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-2.96%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08681"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-4.44%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03449"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.89%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09647"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.71%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001405"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.01%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006355"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3.85%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.578"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-2.66%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3433"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.67%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1290"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-4.38%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.050"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "8.62%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2526"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.31%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04360"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.12%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001238"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.63%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004653"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-4.42%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "176.37%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02205"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "3.20%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05053"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-2.97%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007616"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.60%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009983"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.42%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1373"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.38%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002048"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.83%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008866"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-10.36%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006659"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "4.03%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000757"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.97%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003025"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "13.06%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001206"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "20.64%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002120"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.97%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002019"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.97%"
